$d = $word.ActiveDocument

# 1. Add two trailing spaces after the existing sentence in paragraph 1.
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Microsoft word document.  ", 2)

# 2. Append the red parenthetical note as new runs at the end of paragraph 1.
$para1 = $d.Paragraphs(1)
$endRange = $para1.Range
$endRange.SetRange($endRange.End - 1, $endRange.End - 1)

$r1 = $endRange.InsertAfter("(This is a change " + [char]8211 + " Ve")
$endRange.Font.Color = 255

$endRange.Collapse(0)
$r2 = $endRange.InsertAfter("rsion for main branch")
$endRange.Font.Color = 255

$endRange.Collapse(0)
$r3 = $endRange.InsertAfter(")")
$endRange.Font.Color = 255
